$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(4, 1).Value = -21.54599999999999
$ws.Cells.Item(6, 1).Value = -22.683
$ws.Cells.Item(7, 1).Value = -19.85769999999999
$ws.Cells.Item(7, 3).Value = -12.0729
$ws.Cells.Item(12, 3).Value = -11.00139999999999
$ws.Cells.Item(12, 5).Value = 17.34070000000003
$ws.Cells.Item(13, 4).Value = -8.641999999999998
$ws.Cells.Item(14, 4).Value = -8.008800000000003
$ws.Cells.Item(15, 3).Value = -14.74659999999999
$ws.Cells.Item(16, 1).Value = -21.56819999999999
$ws.Cells.Item(16, 4).Value = -9.177300000000004
$ws.Cells.Item(19, 4).Value = -9.035399999999996
$ws.Cells.Item(20, 1).Value = -18.9406
$ws.Cells.Item(20, 3).Value = -11.8466
$ws.Cells.Item(21, 3).Value = -11.8487
$ws.Cells.Item(22, 3).Value = -12.0847
$ws.Cells.Item(22, 4).Value = -8.4579
$ws.Cells.Item(22, 5).Value = 16.7226
$ws.Cells.Item(23, 3).Value = -11.72620000000001
$ws.Cells.Item(28, 1).Value = -22.04839999999999
$ws.Cells.Item(29, 1).Value = -21.76359999999999
$ws.Cells.Item(29, 3).Value = -11.7127
$ws.Cells.Item(29, 5).Value = 17.23710000000001
$ws.Cells.Item(32, 1).Value = -21.2283
$ws.Cells.Item(34, 3).Value = -11.49500000000001
$ws.Cells.Item(34, 5).Value = 17.4145
$ws.Cells.Item(36, 4).Value = -8.736599999999994
$ws.Cells.Item(40, 1).Value = -20.4875
$ws.Cells.Item(42, 3).Value = -12.2354
$ws.Cells.Item(43, 3).Value = -12.57309999999999
$ws.Cells.Item(43, 5).Value = 17.31660000000002
$ws.Cells.Item(44, 3).Value = -13.93659999999999
$ws.Cells.Item(45, 3).Value = -13.55929999999999
$ws.Cells.Item(46, 1).Value = -21.95670000000001
$ws.Cells.Item(46, 3).Value = -13.8597
$ws.Cells.Item(46, 4).Value = -8.521800000000001
$ws.Cells.Item(48, 5).Value = 17.52630000000001
$ws.Cells.Item(50, 3).Value = -13.86629999999999
$ws.Cells.Item(50, 4).Value = -8.011699999999998
$ws.Cells.Item(51, 1).Value = -22.1325
$ws.Cells.Item(51, 3).Value = -13.02539999999999
$ws.Cells.Item(52, 1).Value = -22.07529999999999
$ws.Cells.Item(57, 1).Value = -22.69700000000001
$ws.Cells.Item(59, 1).Value = -21.9738
$ws.Cells.Item(60, 5).Value = 15.74160000000001
$ws.Cells.Item(62, 1).Value = -21.996
$ws.Cells.Item(66, 1).Value = -22.13530000000001
$ws.Cells.Item(66, 3).Value = -13.52969999999999
$ws.Cells.Item(67, 3).Value = -11.0073
$ws.Cells.Item(68, 5).Value = 17.73070000000002
$ws.Cells.Item(70, 5).Value = 18.50040000000002
$ws.Cells.Item(73, 1).Value = -20.96419999999999
$ws.Cells.Item(73, 5).Value = 17.4019
$ws.Cells.Item(74, 1).Value = -22.15919999999999
$ws.Cells.Item(79, 3).Value = -11.02940000000001
$ws.Cells.Item(84, 3).Value = -13.32349999999998
$ws.Cells.Item(87, 5).Value = 16.36789999999999
$ws.Cells.Item(92, 1).Value = -21.90469999999999
$ws.Cells.Item(92, 3).Value = -11.7125
$ws.Cells.Item(92, 5).Value = 18.07440000000001
$ws.Cells.Item(95, 4).Value = -8.062900000000003
$ws.Cells.Item(97, 3).Value = -11.20410000000001
$ws.Cells.Item(97, 4).Value = -8.146399999999995
$ws.Cells.Item(100, 1).Value = -22.18229999999999
$ws.Cells.Item(101, 5).Value = 16.7677
